# Auto-generated Excel COM-interop script applying "Horarios actualizados Linea 141 - 879" update
$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 04:52:24"
$ws1.Range("A3").Value = "Total filas: 36"

$rows1 = @(
    ,@("03:52:04", "04:01", "81_EL PELIGRO", 9, "LP1912")
    ,@("04:32:18", "04:33", "15_ABASTO", 1, "LP1912")
    ,@("04:44:46", "04:46", "215_EL PELIGRO", 2, "LP1912")
    ,@("03:52:04", "04:46", "215A_EL PATO", 54, "LP1912")
    ,@("04:44:46", "04:46", "15_ABASTO", 2, "LP1912")
    ,@("04:32:18", "04:47", "215_EL PELIGRO", 15, "LP1912")
    ,@("04:52:24", "04:53", "11_ETCHEVERRY", 1, "LP1912")
    ,@("04:52:24", "04:54", "15_ABASTO", 2, "LP1912")
    ,@("04:13:31", "05:11", "17_ROMERO", 58, "LP1912")
    ,@("03:52:04", "05:16", "17_ROMERO", 84, "LP1912")
    ,@("04:52:24", "05:22", "23_HERNANDEZ", 30, "LP1912")
    ,@("04:44:46", "05:31", "81_EL PELIGRO", 47, "LP1912")
    ,@("04:52:24", "05:32", "81_EL PELIGRO", 40, "LP1912")
    ,@("03:52:04", "05:35", "215B_EL PATO", 103, "LP1912")
    ,@("04:52:24", "05:44", "14_ABASTO", 52, "LP1912")
    ,@("03:52:04", "05:46", "15_ABASTO", 114, "LP1912")
    ,@("04:32:18", "05:47", "14_ABASTO", 75, "LP1912")
    ,@("04:13:31", "05:50", "14_ABASTO", 97, "LP1912")
    ,@("04:44:46", "05:51", "17_ROMERO", 67, "LP1912")
    ,@("04:52:24", "05:52", "17_ROMERO", 60, "LP1912")
    ,@("04:44:46", "06:00", "16_SANTA ANA", 76, "LP1912")
    ,@("04:52:24", "06:01", "16_SANTA ANA", 69, "LP1912")
    ,@("04:44:46", "06:03", "10_OLMOS", 79, "LP1912")
    ,@("04:52:24", "06:04", "10_OLMOS", 72, "LP1912")
    ,@("04:44:46", "06:10", "215A_EL PATO", 86, "LP1912")
    ,@("04:52:24", "06:11", "215A_EL PATO", 79, "LP1912")
    ,@("04:32:18", "06:15", "17_ROMERO", 103, "LP1912")
    ,@("04:52:24", "06:24", "11_ETCHEVERRY", 92, "LP1912")
    ,@("04:52:24", "06:27", "23_HERNANDEZ", 95, "LP1912")
    ,@("04:44:46", "06:28", "17_ROMERO", 104, "LP1912")
    ,@("04:44:46", "06:30", "16_SANTA ANA", 106, "LP1912")
    ,@("04:52:24", "06:31", "17X38_ROMERO", 99, "LP1912")
    ,@("04:52:24", "06:31", "16_SANTA ANA", 99, "LP1912")
    ,@("04:52:24", "06:36", "17_ROMERO", 104, "LP1912")
    ,@("04:52:24", "06:39", "225_C ROCA-H SUR", 107, "LP1912")
    ,@("04:52:24", "06:51", "215A_EL PATO", 119, "LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 04:52:24"
$ws2.Range("A3").Value = "Total filas: 7"

$rows2 = @(
    ,@("04:44:46", "04:46", "215_EL PELIGRO", 2, "LP1912")
    ,@("03:52:04", "04:46", "215A_EL PATO", 54, "LP1912")
    ,@("04:32:18", "04:47", "215_EL PELIGRO", 15, "LP1912")
    ,@("03:52:04", "05:35", "215B_EL PATO", 103, "LP1912")
    ,@("04:44:46", "06:10", "215A_EL PATO", 86, "LP1912")
    ,@("04:52:24", "06:11", "215A_EL PATO", 79, "LP1912")
    ,@("04:52:24", "06:51", "215A_EL PATO", 119, "LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 04:52:24"

